# Daily automated update of the "剩余" (remaining days) tracking sheet.
# For every data row:
#   - Column D = total days ("总天")
#   - Column E = remaining days ("剩余")
#   - Column F = start date ("开始时间"), format yyyymmdd
#
# Business rule applied by the nightly job:
#   - If the date in column F is not a well-formed 8-digit yyyymmdd value,
#     the row is skipped (treated as bad data, left untouched).
#   - Otherwise, if the remaining-days counter has already hit the lowest
#     tracked value (1), the item is considered replenished/restarted:
#       E is reset back to the total-days value (D) and F is refreshed to
#       the new start date.
#   - Otherwise the remaining-days counter is simply decremented by 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New start date to stamp on rows that just got refreshed/restocked.
$newStartDate = 20251121

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)   # D: total days
    $eCell = $ws.Cells.Item($r, 5)   # E: remaining days
    $fCell = $ws.Cells.Item($r, 6)   # F: start date (yyyymmdd)

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($dVal -eq $null -or $eVal -eq $null) {
        continue
    }

    # Validate the date cell looks like a proper 8-digit yyyymmdd value;
    # malformed entries are left untouched (skipped), just like the
    # original nightly job does.
    $fText = [string]$fVal
    if ($fText.Length -ne 8) {
        continue
    }

    if ($eVal -le 1) {
        # Remaining days bottomed out -> restock: refill to the full
        # total-days amount and restamp the start date.
        $eCell.Value = $dVal
        $fCell.Value = $newStartDate
    } else {
        $eCell.Value = $eVal - 1
    }
}
